$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2087912087912088
$ws.Range("C2").Value = 0.5357142857142857
$ws.Range("J2").Value = 0.01648351648351648
$ws.Range("P2").Value = 0.1593406593406593
$ws.Range("S2").Value = 0.07967032967032966
$ws.Range("C3").Value = 0.02926829268292683
$ws.Range("J3").Value = 0.02439024390243903
$ws.Range("P3").Value = 0.751219512195122
$ws.Range("S3").Value = 0.1951219512195122
$ws.Range("J4").Value = 0.02777777777777778
$ws.Range("P4").Value = 0.6388888888888888
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.07983193277310924
$ws.Range("D6").Value = 0.004201680672268907
$ws.Range("F6").Value = 0.06722689075630252
$ws.Range("J6").Value = 0.2142857142857143
$ws.Range("Q6").Value = 0.180672268907563
$ws.Range("R6").Value = 0.07142857142857142
$ws.Range("S6").Value = 0.3823529411764706
$ws.Range("B7").Value = 0.1017699115044248
$ws.Range("D7").Value = 0.008849557522123894
$ws.Range("E7").Value = 0.004424778761061947
$ws.Range("F7").Value = 0.05309734513274336
$ws.Range("J7").Value = 0.1460176991150443
$ws.Range("O7").Value = 0.02212389380530973
$ws.Range("Q7").Value = 0.1769911504424779
$ws.Range("R7").Value = 0.084070796460177
$ws.Range("S7").Value = 0.4026548672566372
$ws.Range("B8").Value = 0.1178861788617886
$ws.Range("D8").Value = 0.01829268292682927
$ws.Range("E8").Value = 0.002032520325203252
$ws.Range("F8").Value = 0.0426829268292683
$ws.Range("J8").Value = 0.1219512195121951
$ws.Range("O8").Value = 0.01626016260162602
$ws.Range("Q8").Value = 0.1463414634146341
$ws.Range("R8").Value = 0.1117886178861789
$ws.Range("S8").Value = 0.4227642276422764
$ws.Range("B9").Value = 0.1095890410958904
$ws.Range("D9").Value = 0.0136986301369863
$ws.Range("F9").Value = 0.0639269406392694
$ws.Range("J9").Value = 0.1461187214611872
$ws.Range("O9").Value = 0.0091324200913242
$ws.Range("Q9").Value = 0.1735159817351598
$ws.Range("R9").Value = 0.0958904109589041
$ws.Range("S9").Value = 0.3881278538812785
$ws.Range("B10").Value = 0.1184593023255814
$ws.Range("D10").Value = 0.01598837209302326
$ws.Range("E10").Value = 0.0007267441860465116
$ws.Range("F10").Value = 0.0690406976744186
$ws.Range("J10").Value = 0.1075581395348837
$ws.Range("O10").Value = 0.01162790697674419
$ws.Range("Q10").Value = 0.2020348837209302
$ws.Range("R10").Value = 0.08502906976744186
$ws.Range("S10").Value = 0.3895348837209303
$ws.Range("G11").Value = 0.150887573964497
$ws.Range("J11").Value = 0.07988165680473373
$ws.Range("K11").Value = 0.2159763313609467
$ws.Range("L11").Value = 0.5384615384615384
$ws.Range("S11").Value = 0.01479289940828402
$ws.Range("G12").Value = 0.7553191489361702
$ws.Range("J12").Value = 0.2127659574468085
$ws.Range("L12").Value = 0.01595744680851064
$ws.Range("S12").Value = 0.01595744680851064
$ws.Range("G13").Value = 0.7037037037037037
$ws.Range("J13").Value = 0.2592592592592592
$ws.Range("S13").Value = 0.03703703703703703
$ws.Range("F15").Value = 0.02752293577981652
$ws.Range("H15").Value = 0.2247706422018349
$ws.Range("I15").Value = 0.06880733944954129
$ws.Range("J15").Value = 0.3256880733944954
$ws.Range("K15").Value = 0.07339449541284404
$ws.Range("M15").Value = 0.01834862385321101
$ws.Range("O15").Value = 0.04587155963302753
$ws.Range("S15").Value = 0.2155963302752294
$ws.Range("F16").Value = 0.01327433628318584
$ws.Range("H16").Value = 0.163716814159292
$ws.Range("I16").Value = 0.09734513274336283
$ws.Range("J16").Value = 0.4070796460176991
$ws.Range("K16").Value = 0.1150442477876106
$ws.Range("M16").Value = 0.02654867256637168
$ws.Range("N16").Value = 0.004424778761061947
$ws.Range("O16").Value = 0.03539823008849557
$ws.Range("S16").Value = 0.1371681415929203
$ws.Range("F17").Value = 0.01495726495726496
$ws.Range("H17").Value = 0.1773504273504274
$ws.Range("I17").Value = 0.1004273504273504
$ws.Range("J17").Value = 0.3931623931623932
$ws.Range("K17").Value = 0.08547008547008547
$ws.Range("M17").Value = 0.01923076923076923
$ws.Range("O17").Value = 0.06623931623931624
$ws.Range("S17").Value = 0.1431623931623932
$ws.Range("F18").Value = 0.01777777777777778
$ws.Range("H18").Value = 0.1511111111111111
$ws.Range("I18").Value = 0.09333333333333334
$ws.Range("J18").Value = 0.4266666666666667
$ws.Range("K18").Value = 0.08444444444444445
$ws.Range("M18").Value = 0.008888888888888889
$ws.Range("N18").Value = 0.004444444444444444
$ws.Range("O18").Value = 0.07111111111111111
$ws.Range("S18").Value = 0.1422222222222222
$ws.Range("F19").Value = 0.01644962302947224
$ws.Range("H19").Value = 0.2001370801919123
$ws.Range("I19").Value = 0.07950651130911583
$ws.Range("J19").Value = 0.3625771076079506
$ws.Range("K19").Value = 0.111720356408499
$ws.Range("M19").Value = 0.02604523646333105
$ws.Range("O19").Value = 0.06922549691569568
$ws.Range("S19").Value = 0.1343385880740233
